$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.750.32"
$ws.Range("E2").Value = "  -1.90%  "
$ws.Range("D3").Value = "'1.938.93"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Value = "'242.07"
$ws.Range("E5").Value = "  -2.60%  "
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("D7").Value = "'0.4876"
$ws.Range("E7").Value = "  -0.49%  "
$ws.Range("D8").Value = "'0.2921"
$ws.Range("E8").Value = "  -2.31%  "
$ws.Range("D9").Value = "'0.06865"
$ws.Range("E9").Value = "  -0.21%  "
$ws.Range("D10").Value = "'19.46"
$ws.Range("E10").Value = "  +0.68%  "
$ws.Range("D11").Value = "'105.12"
$ws.Range("E11").Value = "  -2.06%  "
$ws.Range("D12").Value = "'1.980.85"
$ws.Range("E12").Value = "  +1.91%  "
$ws.Range("D13").Value = "'0.07751"
$ws.Range("E13").Value = "  -0.41%  "
$ws.Range("D14").Value = "'5.308"
$ws.Range("E14").Value = "  -3.14%  "
$ws.Range("D15").Value = "'0.6955"
$ws.Range("E15").Value = "  -3.14%  "
$ws.Range("D16").Value = "'274.88"
$ws.Range("E16").Value = "  -5.01%  "
$ws.Range("D17").Value = "'30.766.98"
$ws.Range("E17").Value = "  -1.83%  "
$ws.Range("D18").Value = "'0.000007697"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").Value = "'13.12"
$ws.Range("E19").Value = "  -1.82%  "
$ws.Range("D20").Value = "'2.207.39"
$ws.Range("E20").Value = "  +0.30%  "
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("D22").Value = "'5.440"
$ws.Range("E22").Value = "  -3.89%  "
$ws.Range("D23").Value = "'1.004"
$ws.Range("E23").Value = "  +0.39%  "
$ws.Range("D24").Value = "'6.451"
$ws.Range("E24").Value = "  -3.33%  "
$ws.Range("D25").Value = "'9.686"
$ws.Range("E25").Value = "  -3.84%  "
$ws.Range("D26").Value = "'167.50"
$ws.Range("E26").Value = "  -1.19%  "
$ws.Range("E27").Value = "  -2.72%  "
$ws.Range("D28").Value = "'2.159"
$ws.Range("E28").Value = "  -2.14%  "
$ws.Range("D29").Value = "'0.1038"
$ws.Range("E29").Value = "  -3.05%  "
$ws.Range("D30").Value = "'1.390"
$ws.Range("E30").Value = "  -3.94%  "
$ws.Range("D31").Value = "'1.550"
$ws.Range("E31").Value = "  -3.24%  "
$ws.Range("D32").Value = "'4.541"
$ws.Range("E32").Value = "  -6.82%  "
$ws.Range("D33").Value = "'4.344"
$ws.Range("E33").Value = "  -4.33%  "
$ws.Range("D34").Value = "'0.04839"
$ws.Range("E34").Value = "  -5.04%  "
$ws.Range("D35").Value = "'0.7463"
$ws.Range("E35").Value = "  -3.70%  "
$ws.Range("D36").Value = "'1.153"
$ws.Range("E36").Value = "  -2.23%  "
$ws.Range("D37").Value = "'2.729"
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").Value = "'0.01985"
$ws.Range("E38").Value = "  -3.85%  "
$ws.Range("E39").Value = "  -2.00%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "'77.42"
$ws.Range("E40").Value = "  +4.83%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").Value = "'6.433"
$ws.Range("E41").Value = "  -0.08%  "
$ws.Range("D42").Value = "'2.079"
$ws.Range("E42").Value = "  -3.77%  "
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("D44").Value = "'108.02"
$ws.Range("E44").Value = "  -2.04%  "
$ws.Range("D45").Value = "'0.4398"
$ws.Range("E45").Value = "  -2.43%  "
$ws.Range("D46").Value = "'0.9986"
$ws.Range("E46").Value = "  -0.28%  "
$ws.Range("E47").Value = "  +2.48%  "
$ws.Range("D48").Value = "'1.002.05"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'9.178"
$ws.Range("E50").Value = "  -2.98%  "
$ws.Range("B51").Value = "Elrond"
$ws.Range("C51").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D51").Value = "'35.68"
$ws.Range("E51").Value = "  -1.35%  "
